# Auto-generated edit script applying the Kujata_Profits.xlsx diff
# Updates currentAveragePrice / NQ / HQ / Leve price / profit columns (H-N)
# across several worksheets (rows correspond to specific Leve entries).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 2247.577
$ws.Range("I28").Value = 2084.6
$ws.Range("K28").Value = 2084.6
$ws.Range("M28").Value = -1599.6
# Row 40
$ws.Range("H40").Value = 2478.5715
$ws.Range("I40").Value = 2336.3635
$ws.Range("K40").Value = 2336.3635
$ws.Range("M40").Value = -2161.3635
# Row 62
$ws.Range("H62").Value = 22224984
$ws.Range("I62").Value = 27779728
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 27779728
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -27779104
$ws.Range("N62").Value = -7248
# Row 65
$ws.Range("H65").Value = 22224984
$ws.Range("I65").Value = 27779728
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 138898640
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -138895520
$ws.Range("N65").Value = -36240
# Row 100
$ws.Range("H100").Value = 1070.8823
$ws.Range("I100").Value = 713.2143
$ws.Range("K100").Value = 713.2143
$ws.Range("M100").Value = -172.2143
# Row 125
$ws.Range("H125").Value = 3187.5
$ws.Range("I125").Value = 2311.3333
$ws.Range("J125").Value = 3479.5557
$ws.Range("K125").Value = 20801.9997
$ws.Range("L125").Value = 31316.0013
$ws.Range("M125").Value = -18341.9997
$ws.Range("N125").Value = -36236.0013

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4278.83
$ws.Range("I32").Value = 3854.4639
$ws.Range("J32").Value = 18000
$ws.Range("K32").Value = 3854.4639
$ws.Range("L32").Value = 18000
$ws.Range("M32").Value = -3567.4639
$ws.Range("N32").Value = -18574
# Row 74
$ws.Range("H74").Value = 3001.6191
$ws.Range("I74").Value = 2440.3
$ws.Range("K74").Value = 2440.3
$ws.Range("M74").Value = -1566.3
# Row 77
$ws.Range("H77").Value = 3001.6191
$ws.Range("I77").Value = 2440.3
$ws.Range("K77").Value = 12201.5
$ws.Range("M77").Value = -7833.5
# Row 97
$ws.Range("H97").Value = 625.2222
$ws.Range("I97").Value = 659.5714
$ws.Range("J97").Value = 505
$ws.Range("K97").Value = 659.5714
$ws.Range("L97").Value = 505
$ws.Range("M97").Value = -163.5714
$ws.Range("N97").Value = -1497
# Row 110
$ws.Range("H110").Value = 1973.4286
$ws.Range("I110").Value = 1150.25
$ws.Range("K110").Value = 1150.25
$ws.Range("M110").Value = 894.75
# Row 122
$ws.Range("H122").Value = 3443.7
$ws.Range("I122").Value = 3427.875
$ws.Range("J122").Value = 3507
$ws.Range("K122").Value = 10283.625
$ws.Range("L122").Value = 10521
$ws.Range("M122").Value = -7833.625
$ws.Range("N122").Value = -15421

$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 10000
$ws.Range("I26").Value = 10000
$ws.Range("K26").Value = 10000
$ws.Range("M26").Value = -9708
# Row 86
$ws.Range("H86").Value = 2055.8975
$ws.Range("I86").Value = 1883.7188
$ws.Range("J86").Value = 2843
$ws.Range("K86").Value = 1883.7188
$ws.Range("L86").Value = 2843
$ws.Range("M86").Value = -760.7188000000001
$ws.Range("N86").Value = -5089
# Row 89
$ws.Range("H89").Value = 2055.8975
$ws.Range("I89").Value = 1883.7188
$ws.Range("J89").Value = 2843
$ws.Range("K89").Value = 9418.594000000001
$ws.Range("L89").Value = 14215
$ws.Range("M89").Value = -3802.594000000001
$ws.Range("N89").Value = -25447
# Row 94
$ws.Range("H94").Value = 8621044
$ws.Range("I94").Value = 9615756
$ws.Range("K94").Value = 9615756
$ws.Range("M94").Value = -9615305
# Row 99
$ws.Range("H99").Value = 66667756
$ws.Range("I99").Value = 76924060
$ws.Range("K99").Value = 76924060
$ws.Range("M99").Value = -76922562
# Row 140
$ws.Range("H140").Value = 29680
$ws.Range("J140").Value = 29680
$ws.Range("L140").Value = 29680
$ws.Range("N140").Value = -40040

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 58824690
$ws.Range("I16").Value = 76924264
$ws.Range("K16").Value = 76924264
$ws.Range("M16").Value = -76923977
# Row 22
$ws.Range("H22").Value = 233667.33
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 350351
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 350351
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = -351051
# Row 31
$ws.Range("H31").Value = 1190.0159
$ws.Range("I31").Value = 1060.8077
$ws.Range("J31").Value = 1800.8182
$ws.Range("K31").Value = 1060.8077
$ws.Range("L31").Value = 1800.8182
$ws.Range("M31").Value = -765.8077000000001
$ws.Range("N31").Value = -2390.8182
# Row 34
$ws.Range("H34").Value = 1190.0159
$ws.Range("I34").Value = 1060.8077
$ws.Range("J34").Value = 1800.8182
$ws.Range("K34").Value = 1060.8077
$ws.Range("L34").Value = 1800.8182
$ws.Range("M34").Value = -858.8077000000001
$ws.Range("N34").Value = -2204.8182
# Row 58
$ws.Range("H58").Value = 5193.0967
$ws.Range("I58").Value = 1305.1111
$ws.Range("J58").Value = 10576.462
$ws.Range("K58").Value = 1305.1111
$ws.Range("L58").Value = 10576.462
$ws.Range("M58").Value = -1102.1111
$ws.Range("N58").Value = -10982.462
# Row 105
$ws.Range("H105").Value = 594.5
$ws.Range("I105").Value = 594.5
$ws.Range("K105").Value = 594.5
$ws.Range("M105").Value = 1152.5
# Row 107
$ws.Range("H107").Value = 852
$ws.Range("I107").Value = 455.6154
$ws.Range("J107").Value = 1882.6
$ws.Range("K107").Value = 455.6154
$ws.Range("L107").Value = 1882.6
$ws.Range("M107").Value = 1464.3846
$ws.Range("N107").Value = -5722.6
# Row 113
$ws.Range("H113").Value = 58824690
$ws.Range("I113").Value = 76924264
$ws.Range("K113").Value = 76924264
$ws.Range("M113").Value = -76922094
# Row 122
$ws.Range("H122").Value = 682.3182
$ws.Range("J122").Value = 510.5
$ws.Range("L122").Value = 1531.5
$ws.Range("N122").Value = -6431.5
# Row 136
$ws.Range("I136").Value = 1305.1111
$ws.Range("J136").Value = 10576.462
$ws.Range("K136").Value = 3915.3333
$ws.Range("L136").Value = 31729.386
$ws.Range("M136").Value = -1365.3333
$ws.Range("N136").Value = -36829.386

$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 121.85714
$ws.Range("I38").Value = 125.5
$ws.Range("K38").Value = 376.5
$ws.Range("M38").Value = -29.5
# Row 107
$ws.Range("H107").Value = 3784.5334
$ws.Range("J107").Value = 4755.5654
$ws.Range("L107").Value = 14266.6962
$ws.Range("N107").Value = -18106.6962
# Row 113
$ws.Range("H113").Value = 670.2045000000001
$ws.Range("I113").Value = 582.3570999999999
$ws.Range("J113").Value = 711.2
$ws.Range("K113").Value = 1747.0713
$ws.Range("L113").Value = 2133.6
$ws.Range("M113").Value = 422.9287000000002
$ws.Range("N113").Value = -6473.6
# Row 122
$ws.Range("H122").Value = 1308.5
$ws.Range("I122").Value = 606.6
$ws.Range("J122").Value = 2010.4
$ws.Range("K122").Value = 5459.400000000001
$ws.Range("L122").Value = 18093.6
$ws.Range("M122").Value = -3009.400000000001
$ws.Range("N122").Value = -22993.6
# Row 131
$ws.Range("H131").Value = 725.3
$ws.Range("I131").Value = 453.27274
$ws.Range("J131").Value = 758.9213
$ws.Range("K131").Value = 1359.81822
$ws.Range("L131").Value = 2276.7639
$ws.Range("M131").Value = 3680.18178
$ws.Range("N131").Value = -12356.7639

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
# Row 118
$ws.Range("H118").Value = 20271.25
$ws.Range("J118").Value = 20271.25
$ws.Range("L118").Value = 20271.25
$ws.Range("N118").Value = -23585.25

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 764.8182
$ws.Range("I22").Value = 676.375
$ws.Range("J22").Value = 1000.6667
$ws.Range("K22").Value = 676.375
$ws.Range("L22").Value = 1000.6667
$ws.Range("M22").Value = -381.375
$ws.Range("N22").Value = -1590.6667
# Row 27
$ws.Range("H27").Value = 764.8182
$ws.Range("I27").Value = 676.375
$ws.Range("J27").Value = 1000.6667
$ws.Range("K27").Value = 676.375
$ws.Range("L27").Value = 1000.6667
$ws.Range("M27").Value = -569.375
$ws.Range("N27").Value = -1214.6667
# Row 32
$ws.Range("H32").Value = 5416.6665
$ws.Range("I32").Value = 4125
$ws.Range("K32").Value = 4125
$ws.Range("M32").Value = -3808
# Row 68
$ws.Range("H68").Value = 1288.4615
$ws.Range("I68").Value = 1288.4615
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1288.4615
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -539.4614999999999
$ws.Range("N68").ClearContents()
# Row 71
$ws.Range("H71").Value = 1288.4615
$ws.Range("I71").Value = 1288.4615
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 6442.307499999999
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -2698.307499999999
$ws.Range("N71").ClearContents()
# Row 100
$ws.Range("H100").Value = 1300
$ws.Range("I100").Value = 1000
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 1000
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -459
$ws.Range("N100").Value = -2582
# Row 111
$ws.Range("H111").Value = 65000
$ws.Range("J111").Value = 65000
$ws.Range("L111").Value = 65000
$ws.Range("N111").Value = -73180

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 166669000
$ws.Range("I62").Value = 166669000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 166669000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -166668376
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 166669000
$ws.Range("I65").Value = 166669000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 833345000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -833341880
$ws.Range("N65").ClearContents()
# Row 107
$ws.Range("H107").Value = 463.77274
$ws.Range("J107").Value = 538.9167
$ws.Range("L107").Value = 1616.7501
$ws.Range("N107").Value = -5456.7501
# Row 122
$ws.Range("H122").Value = 11906626
$ws.Range("I122").Value = 13159834
$ws.Range("K122").Value = 39479502
$ws.Range("M122").Value = -39477052
